$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "51.510.71"
Set-TextValue $ws.Range("E2") "  +4.82%  "
Set-TextValue $ws.Range("D3") "2.741.89"
Set-TextValue $ws.Range("E3") "  +4.26%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.08%  "
Set-TextValue $ws.Range("D5") "115.80"
Set-TextValue $ws.Range("E5") "  +4.03%  "
Set-TextValue $ws.Range("D6") "331.99"
Set-TextValue $ws.Range("E6") "  +2.89%  "
Set-TextValue $ws.Range("D7") "0.537"
Set-TextValue $ws.Range("E7") "  +2.22%  "
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.06%  "
Set-TextValue $ws.Range("D9") "0.571"
Set-TextValue $ws.Range("E9") "  +4.97%  "
Set-TextValue $ws.Range("D10") "41.51"
Set-TextValue $ws.Range("E10") "  +4.44%  "
Set-TextValue $ws.Range("D11") "0.0854"
Set-TextValue $ws.Range("E11") "  +5.35%  "
Set-TextValue $ws.Range("D12") "20.10"
Set-TextValue $ws.Range("E12") "  +1.37%  "
Set-TextValue $ws.Range("D13") "0.129"
Set-TextValue $ws.Range("E13") "  +2.84%  "
Set-TextValue $ws.Range("D14") "7.60"
Set-TextValue $ws.Range("E14") "  +4.83%  "
Set-TextValue $ws.Range("D15") "3.162.19"
Set-TextValue $ws.Range("E15") "  +3.86%  "
Set-TextValue $ws.Range("D16") "2.756.44"
Set-TextValue $ws.Range("E16") "  +4.91%  "
Set-TextValue $ws.Range("D17") "0.880"
Set-TextValue $ws.Range("E17") "  +2.28%  "
Set-TextValue $ws.Range("D18") "51.426.43"
Set-TextValue $ws.Range("E18") "  +4.66%  "
Set-TextValue $ws.Range("D19") "3.21"
Set-TextValue $ws.Range("E19") "  +7.80%  "
Set-TextValue $ws.Range("D20") "13.39"
Set-TextValue $ws.Range("E20") "  +3.97%  "
Set-TextValue $ws.Range("D21") "6.83"
Set-TextValue $ws.Range("E21") "  +2.03%  "
Set-TextValue $ws.Range("D22") "0.0₃0974"
Set-TextValue $ws.Range("E22") "  +3.02%  "
Set-TextValue $ws.Range("D23") "278.18"
Set-TextValue $ws.Range("E23") "  +3.24%  "
Set-TextValue $ws.Range("D24") "69.31"
Set-TextValue $ws.Range("E24") "  +1.17%  "
Set-TextValue $ws.Range("D25") "2.64"
Set-TextValue $ws.Range("E25") "  +4.29%  "
Set-TextValue $ws.Range("D26") "26.74"
Set-TextValue $ws.Range("E26") "  +2.46%  "
Set-TextValue $ws.Range("D27") "1.00"
Set-TextValue $ws.Range("E27") "  -0.08%  "
Set-TextValue $ws.Range("D28") "10.15"
Set-TextValue $ws.Range("E28") "  +1.18%  "
Set-TextValue $ws.Range("D29") "2.22"
Set-TextValue $ws.Range("E29") "  -0.14%  "
Set-TextValue $ws.Range("D30") "0.141"
Set-TextValue $ws.Range("E30") "  +1.89%  "
Set-TextValue $ws.Range("D31") "34.99"
Set-TextValue $ws.Range("E31") "  -0.23%  "
Set-TextValue $ws.Range("D32") "50.01"
Set-TextValue $ws.Range("E32") "  +1.02%  "
Set-TextValue $ws.Range("D33") "5.53"
Set-TextValue $ws.Range("E33") "  +0.78%  "
Set-TextValue $ws.Range("D34") "0.0819"
Set-TextValue $ws.Range("E34") "  +2.65%  "
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  -0.11%  "
Set-TextValue $ws.Range("D36") "19.07"
Set-TextValue $ws.Range("E36") "  +0.29%  "
Set-TextValue $ws.Range("D37") "5.00"
Set-TextValue $ws.Range("E37") "  +0.39%  "
Set-TextValue $ws.Range("D38") "2.08"
Set-TextValue $ws.Range("E38") "  +1.95%  "
Set-TextValue $ws.Range("D39") "3.20"
Set-TextValue $ws.Range("E39") "  +2.18%  "
Set-TextValue $ws.Range("D40") "127.64"
Set-TextValue $ws.Range("E40") "  +0.42%  "
Set-TextValue $ws.Range("D41") "23.02"
Set-TextValue $ws.Range("E41") "  +4.36%  "
Set-TextValue $ws.Range("D42") "2.30"
Set-TextValue $ws.Range("E42") "  +7.84%  "
Set-TextValue $ws.Range("D43") "0.113"
Set-TextValue $ws.Range("E43") "  +2.54%  "
Set-TextValue $ws.Range("D44") "0.0343"
Set-TextValue $ws.Range("E44") "  +8.44%  "
Set-TextValue $ws.Range("D45") "2.41"
Set-TextValue $ws.Range("E45") "  +12.04%  "
Set-TextValue $ws.Range("D46") "2.086.00"
Set-TextValue $ws.Range("E46") "  +1.03%  "
Set-TextValue $ws.Range("D47") "3.31"
Set-TextValue $ws.Range("E47") "  +1.95%  "
Set-TextValue $ws.Range("D48") "2.22"
Set-TextValue $ws.Range("E48") "  +3.60%  "
Set-TextValue $ws.Range("D49") "5.52"
Set-TextValue $ws.Range("E49") "  +6.32%  "
Set-TextValue $ws.Range("D50") "8.92"
Set-TextValue $ws.Range("E50") "  +0.57%  "
Set-TextValue $ws.Range("D51") "59.87"
Set-TextValue $ws.Range("E51") "  +2.22%  "
